# Update the "Services" sheet (sheet2) column layout and shared strings,
# then make "Services" the active tab (previously "Blogs" was active).
#
# Old header row (Services!A1:E1): title | subheading | content | description | image_urls
# New header row (Services!A1:I1): title | metaTitle | slug | metaDescription | subheading |
#                                   category | description | content | image_urls
#
# Existing values are relocated (cut/paste) rather than overwritten so that the
# shared-string table keeps its original entries/order; only the three brand new
# strings (metaTitle, metaDescription, category) are appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Services")

# Move existing header cells to their new positions (process right-to-left /
# destination-safe order so a cell isn't overwritten before it is relocated).
$ws.Range("E1").Cut($ws.Range("I1"))   # image_urls: E1 -> I1
$ws.Range("D1").Cut($ws.Range("G1"))   # description: D1 -> G1
$ws.Range("C1").Cut($ws.Range("H1"))   # content:     C1 -> H1
$ws.Range("B1").Cut($ws.Range("E1"))   # subheading:  B1 -> E1

# Fill in the newly opened columns.
$ws.Range("B1").Value = "metaTitle"
$ws.Range("C1").Value = "slug"
$ws.Range("D1").Value = "metaDescription"
$ws.Range("F1").Value = "category"

# Make Services the active sheet/tab, with I1 selected.
$ws.Activate()
$ws.Range("I1").Select()
